# Rename the VEO-authorization / coordinator-verification fields to the
# clearer "_as" names and switch their choice list from true/false/unsure
# to valid/invalid/unsure, updating labels accordingly.

$wb = $excel.ActiveWorkbook

# --- survey sheet -----------------------------------------------------
$survey = $wb.Worksheets.Item("survey")

# Row 43: "has_been_verified_by_coordinator" -> "has_been_verified_by_coordinator_as"
$survey.Range("E43").Value = "valid_invalid_unsure"
$survey.Range("F43").Value = "has_been_verified_by_coordinator_as"
$survey.Range("G43").Value = "Coordinator verified this business as:"

# Row 45: "has_been_authorized_by_veo" -> "has_been_authorized_by_veo_as"
$survey.Range("E45").Value = "valid_invalid_unsure"
$survey.Range("F45").Value = "has_been_authorized_by_veo_as"
$survey.Range("G45").Value = "VEO authorized this business as:"

# --- choices sheet ------------------------------------------------------
$choices = $wb.Worksheets.Item("choices")

# true_false_unsure choice list -> valid_invalid_unsure, with new labels
$choices.Range("A11").Value = "valid_invalid_unsure"
$choices.Range("C11").Value = "Invalid"
$choices.Range("D11").Value = "Sio halali"

$choices.Range("A12").Value = "valid_invalid_unsure"
$choices.Range("C12").Value = "Valid"
$choices.Range("D12").Value = "Halali"

$choices.Range("A13").Value = "valid_invalid_unsure"

# --- model sheet ----------------------------------------------------
$model = $wb.Worksheets.Item("model")

$model.Range("B24").Value = "has_been_verified_by_coordinator_as"
$model.Range("B26").Value = "has_been_authorized_by_veo_as"

# --- view bookkeeping ---------------------------------------------------
# workbook-level active tab marker was cleared
$wb.Windows.Item(1).Activate()

# survey sheet: tab selected, scrolled + new active cell
$survey.Select()
$survey.Application.ActiveWindow.ScrollRow = 27
$survey.Range("G46").Select()

# prompt_types sheet: no longer the selected tab, and selection/scroll reset
$promptTypes = $wb.Worksheets.Item("prompt_types")
$promptTypes.Range("B26").Select()

# model sheet: scroll position reset and new active cell
$model.Range("G13").Select()
